$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 37500
$ws.Range("J3").Value = 37500
$ws.Range("L3").Value = 37500
$ws.Range("N3").Value = -37728

$ws.Range("H18").Value = 1498
$ws.Range("I18").Value = 1498
$ws.Range("K18").Value = 1498
$ws.Range("M18").Value = -1214

$ws.Range("H63").Value = 80000
$ws.Range("J63").Value = 80000
$ws.Range("L63").Value = 80000
$ws.Range("N63").Value = -81248

$ws.Range("H66").Value = 80000
$ws.Range("J66").Value = 80000
$ws.Range("L66").Value = 240000
$ws.Range("N66").Value = -246240

$ws.Range("H69").Value = 6188.8335
$ws.Range("J69").Value = 8283.25
$ws.Range("L69").Value = 24849.75
$ws.Range("N69").Value = -26597.75

$ws.Range("H72").Value = 6188.8335
$ws.Range("J72").Value = 8283.25
$ws.Range("L72").Value = 74549.25
$ws.Range("N72").Value = -83285.25

$ws.Range("H80").Value = 3045.182
$ws.Range("I80").Value = 3221.1428
$ws.Range("J80").Value = 2737.25
$ws.Range("K80").Value = 9663.428400000001
$ws.Range("L80").Value = 8211.75
$ws.Range("M80").Value = -8665.428400000001
$ws.Range("N80").Value = -10207.75

$ws.Range("H83").Value = 3045.182
$ws.Range("I83").Value = 3221.1428
$ws.Range("J83").Value = 2737.25
$ws.Range("K83").Value = 28990.2852
$ws.Range("L83").Value = 24635.25
$ws.Range("M83").Value = -23998.2852
$ws.Range("N83").Value = -34619.25

$ws.Range("H98").Value = 1856.375
$ws.Range("I98").Value = 1986
$ws.Range("J98").Value = 949
$ws.Range("K98").Value = 1986
$ws.Range("L98").Value = 949
$ws.Range("M98").Value = -488
$ws.Range("N98").Value = -3945

$ws.Range("H102").Value = 37500
$ws.Range("J102").Value = 37500
$ws.Range("L102").Value = 37500
$ws.Range("N102").Value = -43990

$ws.Range("H122").Value = 1856.375
$ws.Range("I122").Value = 1986
$ws.Range("J122").Value = 949
$ws.Range("K122").Value = 5958
$ws.Range("L122").Value = 2847
$ws.Range("M122").Value = -3508
$ws.Range("N122").Value = -7747

$ws.Range("H125").Value = 3333
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 3333
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 29997
$ws.Range("N125").Value = -34917
$ws.Range("M125").ClearContents()

$ws.Range("H137").Value = 7221.278
$ws.Range("I137").Value = 1477.125
$ws.Range("K137").Value = 4431.375
$ws.Range("M137").Value = -1881.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2332.476
$ws.Range("I122").Value = 2373.875
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 7121.625
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -4671.625
$ws.Range("N122").Value = -11500

$ws.Range("H132").Value = 2688.8333
$ws.Range("I132").Value = 2545.5908
$ws.Range("K132").Value = 7636.7724
$ws.Range("M132").Value = -5106.7724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1374.6709
$ws.Range("I134").Value = 1352.0139
$ws.Range("K134").Value = 4056.0417
$ws.Range("M134").Value = -1521.0417

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24009.66
$ws.Range("I31").Value = 33428
$ws.Range("J31").Value = 3917.2
$ws.Range("K31").Value = 33428
$ws.Range("L31").Value = 3917.2
$ws.Range("M31").Value = -33133
$ws.Range("N31").Value = -4507.2

$ws.Range("H34").Value = 24009.66
$ws.Range("I34").Value = 33428
$ws.Range("J34").Value = 3917.2
$ws.Range("K34").Value = 33428
$ws.Range("L34").Value = 3917.2
$ws.Range("M34").Value = -33226
$ws.Range("N34").Value = -4321.2

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H122").Value = 2999.5
$ws.Range("I122").Value = 2999.5
$ws.Range("K122").Value = 8998.5
$ws.Range("M122").Value = -6548.5

$ws.Range("H127").Value = 88369.57000000001
$ws.Range("I127").Value = 49997
$ws.Range("J127").Value = 94765
$ws.Range("K127").Value = 49997
$ws.Range("L127").Value = 94765
$ws.Range("M127").Value = -45037
$ws.Range("N127").Value = -104685

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4629970
$ws.Range("I34").Value = 120.888885
$ws.Range("J34").Value = 18519518
$ws.Range("K34").Value = 362.666655
$ws.Range("L34").Value = 55558554
$ws.Range("M34").Value = -278.666655
$ws.Range("N34").Value = -55558722

$ws.Range("H121").Value = 631.3333
$ws.Range("J121").Value = 999
$ws.Range("L121").Value = 2997
$ws.Range("N121").Value = -5617

$ws.Range("H122").Value = 1712
$ws.Range("I122").Value = 3248.75
$ws.Range("J122").Value = 1153.1818
$ws.Range("K122").Value = 29238.75
$ws.Range("L122").Value = 10378.6362
$ws.Range("M122").Value = -26788.75
$ws.Range("N122").Value = -15278.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 29949.5
$ws.Range("J52").Value = 29949.5
$ws.Range("L52").Value = 29949.5
$ws.Range("N52").Value = -30467.5

$ws.Range("H70").Value = 15208.786
$ws.Range("I70").Value = 5201.5
$ws.Range("K70").Value = 5201.5
$ws.Range("M70").Value = -4931.5

$ws.Range("H73").Value = 15208.786
$ws.Range("I73").Value = 5201.5
$ws.Range("K73").Value = 5201.5
$ws.Range("M73").Value = -4265.5

$ws.Range("H80").Value = 7291.0835
$ws.Range("J80").Value = 8432.666999999999
$ws.Range("L80").Value = 8432.666999999999
$ws.Range("N80").Value = -10428.667

$ws.Range("H83").Value = 7291.0835
$ws.Range("J83").Value = 8432.666999999999
$ws.Range("L83").Value = 42163.335
$ws.Range("N83").Value = -52147.335

$ws.Range("H122").Value = 2110.389
$ws.Range("I122").Value = 2049.1875
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 6147.5625
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = -3697.5625
$ws.Range("N122").Value = -12700

$ws.Range("H126").Value = 22978.846
$ws.Range("I126").Value = 34818.57
$ws.Range("K126").Value = 104455.71
$ws.Range("M126").Value = -101985.71

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13821.889
$ws.Range("I7").Value = 18879.8
$ws.Range("J7").Value = 7499.5
$ws.Range("K7").Value = 18879.8
$ws.Range("L7").Value = 7499.5
$ws.Range("M7").Value = -18767.8
$ws.Range("N7").Value = -7723.5

$ws.Range("H22").Value = 3103.6155
$ws.Range("J22").Value = 4062.25
$ws.Range("L22").Value = 4062.25
$ws.Range("N22").Value = -4652.25

$ws.Range("H27").Value = 3103.6155
$ws.Range("J27").Value = 4062.25
$ws.Range("L27").Value = 4062.25
$ws.Range("N27").Value = -4276.25

$ws.Range("H40").Value = 3798.9333
$ws.Range("I40").Value = 3596.8572
$ws.Range("K40").Value = 3596.8572
$ws.Range("M40").Value = -3460.8572

$ws.Range("H45").Value = 8000
$ws.Range("I45").Value = 8000
$ws.Range("K45").Value = 8000
$ws.Range("M45").Value = -7593

$ws.Range("H48").Value = 11333
$ws.Range("I48").Value = 8599.6
$ws.Range("K48").Value = 8599.6
$ws.Range("M48").Value = -7938.6

$ws.Range("H122").Value = 1118383.5
$ws.Range("J122").Value = 8778.6
$ws.Range("L122").Value = 26335.8
$ws.Range("N122").Value = -31235.8

$ws.Range("H126").Value = 13821.889
$ws.Range("I126").Value = 18879.8
$ws.Range("J126").Value = 7499.5
$ws.Range("K126").Value = 56639.39999999999
$ws.Range("L126").Value = 22498.5
$ws.Range("M126").Value = -54169.39999999999
$ws.Range("N126").Value = -27438.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3508.25
$ws.Range("I122").Value = 3388.7778
$ws.Range("K122").Value = 10166.3334
$ws.Range("M122").Value = -7716.3334

$ws.Range("H126").Value = 8540.9375
$ws.Range("I126").Value = 11273.454
$ws.Range("K126").Value = 33820.362
$ws.Range("M126").Value = -31350.362

$ws.Range("H136").Value = 2756.2856
$ws.Range("I136").Value = 2883.3845
$ws.Range("J136").Value = 2549.75
$ws.Range("K136").Value = 8650.1535
$ws.Range("L136").Value = 7649.25
$ws.Range("M136").Value = -6100.1535
$ws.Range("N136").Value = -12749.25
